# Clean up Agency table and add AgencyType
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "AgencyType" worksheet as the last tab.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "AgencyType"

# Column widths / styling (B narrow "code" column, C wider "description" column)
$newSheet.Columns.Item(2).ColumnWidth = 8.2
$newSheet.Columns.Item(3).ColumnWidth = 23.0

# Header row
$newSheet.Range("A1").Value = "ID"
$newSheet.Range("B1").Value = "Code"
$newSheet.Range("C1").Value = "Description"

# Descriptions first (column C), top to bottom.
$newSheet.Range("C2").Value = "Covered by another agency"
$newSheet.Range("C3").Value = "City"
$newSheet.Range("C4").Value = "County"
$newSheet.Range("C5").Value = "University or college"
$newSheet.Range("C6").Value = "State Police"
$newSheet.Range("C7").Value = "Special Agency"
$newSheet.Range("C8").Value = "Other state agencies"
$newSheet.Range("C9").Value = "Tribal agencies"
$newSheet.Range("C10").Value = "Federal agencies"
$newSheet.Range("C11").Value = "Unknown"

# ---------------------------------------------------------------------
# 2. Update the TOC sheet with the new "AgencyType" row.
# ---------------------------------------------------------------------
$toc = $wb.Worksheets.Item("TOC")
$toc.Range("A30").Value = "AgencyType"
$toc.Range("B30").Value = "AgencyType"
[void]$toc.Range("A31").Select()

# ---------------------------------------------------------------------
# 3. Fill in the ID / Code columns of the AgencyType sheet.
# ---------------------------------------------------------------------
$newSheet.Range("A2").Value = 1
$newSheet.Range("A3").Value = 2
$newSheet.Range("A4").Value = 3
$newSheet.Range("A5").Value = 4
$newSheet.Range("A6").Value = 5
$newSheet.Range("A7").Value = 6
$newSheet.Range("A8").Value = 7
$newSheet.Range("A9").Value = 8
$newSheet.Range("A10").Value = 9
$newSheet.Range("A11").Value = 99

$newSheet.Range("B2").Value = "'0"
$newSheet.Range("B3").Value = "'1"
$newSheet.Range("B4").Value = "'2"
$newSheet.Range("B5").Value = "'3"
$newSheet.Range("B6").Value = "'4"
$newSheet.Range("B7").Value = "'5"
$newSheet.Range("B8").Value = "'6"
$newSheet.Range("B9").Value = "'7"
$newSheet.Range("B10").Value = "'8"
$newSheet.Range("B11").Value = "'9"

# ---------------------------------------------------------------------
# 4. Leave the AgencyType sheet active/selected, matching the authored
#    workbook (it becomes the last-edited / visible tab).
# ---------------------------------------------------------------------
[void]$newSheet.Range("B12").Select()
